$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# The "Experimental" row (B7) was blank; FHIR ValueSet export now always
# emits the boolean "experimental" element, so the sheet gains a literal
# text value "true" here (matches how the rest of the sheet stores
# everything as text, not as a native Excel boolean).
$ws.Range("B7").Formula = "=T(""true"")"
$ws.Range("B7").Copy() | Out-Null
$ws.Range("B7").PasteSpecial(-4163) | Out-Null
$excel.CutCopyMode = $false

# The metadata export timestamp moved forward to the day the valueset was
# regenerated with the new "experimental" element.
$ws.Range("B8").Value = "2023-02-01T09:05:11-06:00"
